$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-6: clear the y_0_forecast (C) and y_1_forecast (E) cells entirely
$ws.Range("E2").ClearContents()

$ws.Range("C3").ClearContents()
$ws.Range("E3").ClearContents()

$ws.Range("C4").ClearContents()
$ws.Range("E4").ClearContents()

$ws.Range("C5").ClearContents()
$ws.Range("E5").ClearContents()

$ws.Range("C6").ClearContents()
$ws.Range("E6").ClearContents()

# Rows 7-19: update forecast values (C = y_0_forecast, E = y_1_forecast)
$ws.Range("C7").Value = 1.133560223479058
$ws.Range("E7").Value = 1.985690391709771

$ws.Range("C8").Value = 3.633318781899142
$ws.Range("E8").Value = 2.715291551682419

$ws.Range("C9").Value = 3.057638025163611
$ws.Range("E9").Value = 2.42782168586293

$ws.Range("C10").Value = 2.319057151538662
$ws.Range("E10").Value = 2.508920621023392

$ws.Range("C11").Value = 2.536029549059826
$ws.Range("E11").Value = 2.546671316138061

$ws.Range("C12").Value = 3.120740332206995
$ws.Range("E12").Value = 2.775533179497169

$ws.Range("C13").Value = 2.891533899000343
$ws.Range("E13").Value = 2.545843589346886

$ws.Range("C14").Value = 2.618329006605924
$ws.Range("E14").Value = 2.671430903007876

$ws.Range("C15").Value = 2.137626121054947
$ws.Range("E15").Value = 2.891950990452763

$ws.Range("C16").Value = 4.951039758187648
$ws.Range("E16").Value = 3.481452844954491

$ws.Range("C17").Value = 2.838865660558509
$ws.Range("E17").Value = 2.377254777217375

$ws.Range("C18").Value = 1.625773169906108
$ws.Range("E18").Value = 2.42082970885531

$ws.Range("C19").Value = 2.030491763452114
$ws.Range("E19").Value = 2.559374235215039
